# "Add files via upload" — the author re-uploaded the deposit table after
# changing the interest-rate type for every deposit row from fixed
# ("Фиксна") to variable ("Варијабилна"). All other data (bank, currency,
# amount, term, rate, notes, link) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 2).Value = "Варијабилна"
}

# Match the author's final selection/scroll position in the saved file.
$ws.Range("A30").Select() | Out-Null
